$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: manual edit of G5 (2 -> 1); B5/B6/C6/C7 cascade automatically via formulas ---
$ws.Range("G5").Value = 1

# --- Row 7: layer type switches from Upsampling to ConvTrans, with a new formula
#     and explicit (unstyled) values for D7:H7 ---
$ws.Range("A7").Value = "ConvTrans"
$ws.Range("B7").Formula = "=(C7-1)*E7-2*D7+G7*(F7-1)+H7+1"
$ws.Range("D7:H7").ClearFormats() | Out-Null
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 0

# --- Row 8: manual edit of E8 (1 -> 2); B8 cascades automatically ---
$ws.Range("E8").Value = 2

# --- Row 10: manual edit of E10 (2 -> 3); B10 cascades automatically ---
$ws.Range("E10").Value = 3

# --- Remove rows 11 and 12 entirely (shifts the trailing note row 16 up to row 14) ---
$ws.Rows("11:12").Delete() | Out-Null

# --- Update the active selection to match the saved view state ---
$ws.Range("E10").Select() | Out-Null
